$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Prototipo do Site" task row (row 18) is being removed from the
# backlog - its content (HTML/CSS work) moved over to the web-data-viz
# project. Deleting the entire row shifts everything below it up by one,
# which is exactly what the target workbook shows (the old row 19
# "Padronizar CSS" becomes the new row 18, and the trailing blank rows
# shift from 25-27 to 24-26).
$ws.Rows(18).Delete()

# Clear the leftover ";" placeholder that used to sit in L5.
$ws.Range("L5").ClearContents()
